# Horarios actualizados Línea 141 - 90
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header refresh
$ws1.Range("A2").Value = "Última actualización: 16:33:54"
$ws1.Range("A3").Value = "Total filas: 133"

# Rows 38/39 - swap the "Linea" values
$ws1.Range("C38").Value = "215A_LA PLATA"
$ws1.Range("C39").Value = "14_ABASTO"

# Rows 55/56/57 - rotate Hora_Scrap / Linea / Minutos
$ws1.Range("A55").Value = "13:54:15"
$ws1.Range("C55").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D55").Value = 91

$ws1.Range("A56").Value = "13:35:25"
$ws1.Range("C56").Value = "215C_EL PATO"
$ws1.Range("D56").Value = 110

$ws1.Range("A57").Value = "13:54:15"
$ws1.Range("C57").Value = "11_ETCHEVERRY"
$ws1.Range("D57").Value = 91

# Insert a new scraped row before the old row 108
$ws1.Rows.Item(108).EntireRow.Insert()
$ws1.Range("A108").Value = "16:33:53"
$ws1.Range("B108").Value = "17:13"
$ws1.Range("C108").Value = "10_OLMOS"
$ws1.Range("D108").Value = 40
$ws1.Range("E108").Value = "LP1912"

# Insert a second new scraped row before the (now shifted) old row 127
$ws1.Rows.Item(128).EntireRow.Insert()
$ws1.Range("A128").Value = "16:33:53"
$ws1.Range("B128").Value = "17:47"
$ws1.Range("C128").Value = "27_EL RETIRO"
$ws1.Range("D128").Value = 74
$ws1.Range("E128").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 16:33:54"

# ---------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 16:33:54"
$ws3.Range("A3").Value = "Total filas: 8"

# Insert a new scraped row before the old row 8
$ws3.Rows.Item(8).EntireRow.Insert()
$ws3.Range("A8").Value = "16:33:53"
$ws3.Range("B8").Value = "17:02"
$ws3.Range("C8").Value = "215C_LA PLATA"
$ws3.Range("D8").Value = 29
$ws3.Range("E8").Value = "L6203"

Write-Host "Edit complete"
